$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19, shifting existing rows 19:105 down to 20:106.
$ws.Rows("19:19").Insert()

# Fill the newly inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44635
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 100112052
$ws.Range("G19").Value = "Albahaca"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 90
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 5000
$ws.Range("N19").Value = "$/docena de matas"
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 833
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = "Hortaliza"
